$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 3717.0435
$ws.Range("I69").Value = 3011.8333
$ws.Range("J69").Value = 4486.364
$ws.Range("K69").Value = 9035.499899999999
$ws.Range("L69").Value = 13459.092
$ws.Range("M69").Value = -8161.499899999999
$ws.Range("N69").Value = -15207.092

# Row 72
$ws.Range("H72").Value = 3717.0435
$ws.Range("I72").Value = 3011.8333
$ws.Range("J72").Value = 4486.364
$ws.Range("K72").Value = 27106.4997
$ws.Range("L72").Value = 40377.276
$ws.Range("M72").Value = -22738.4997
$ws.Range("N72").Value = -49113.276

# Row 112
$ws.Range("H112").Value = 1059.2222
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1065.3462
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 3196.0386
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -5412.0386

# Row 116
$ws.Range("H116").Value = 11675.5
$ws.Range("I116").Value = 18192.5
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 18192.5
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = -14750.5
$ws.Range("N116").Value = -8784

# Row 127
$ws.Range("H127").Value = 949.44446
$ws.Range("I127").Value = 893.75
$ws.Range("J127").Value = 994
$ws.Range("K127").Value = 2681.25
$ws.Range("L127").Value = 2982
$ws.Range("M127").Value = 2278.75
$ws.Range("N127").Value = -12902

# Row 129
$ws.Range("H129").Value = 995.1852
$ws.Range("I129").Value = 490
$ws.Range("J129").Value = 1035.6
$ws.Range("K129").Value = 1470
$ws.Range("L129").Value = 3106.8
$ws.Range("M129").Value = 3530
$ws.Range("N129").Value = -13106.8

# Row 140
$ws.Range("H140").Value = 70192.164
$ws.Range("I140").Value = 30000
$ws.Range("J140").Value = 72556.414
$ws.Range("K140").Value = 30000
$ws.Range("L140").Value = 72556.414
$ws.Range("M140").Value = -24820
$ws.Range("N140").Value = -82916.414

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5868.7173
$ws.Range("I32").Value = 2926.861
$ws.Range("K32").Value = 2926.861
$ws.Range("M32").Value = -2639.861

# Row 45
$ws.Range("H45").Value = 3130.2856
$ws.Range("I45").Value = 2304
$ws.Range("J45").Value = 3750
$ws.Range("K45").Value = 2304
$ws.Range("L45").Value = 3750
$ws.Range("M45").Value = -1927
$ws.Range("N45").Value = -4504

# Row 51
$ws.Range("H51").Value = 29900
$ws.Range("J51").Value = 29900
$ws.Range("L51").Value = 29900
$ws.Range("N51").Value = -31412

# Row 52
$ws.Range("H52").Value = 43126.668
$ws.Range("J52").Value = 43126.668
$ws.Range("L52").Value = 43126.668
$ws.Range("N52").Value = -43762.668

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 93
$ws.Range("H93").Value = 63155
$ws.Range("J93").Value = 63155
$ws.Range("L93").Value = 63155
$ws.Range("N93").Value = -68147

# Row 119
$ws.Range("H119").Value = 39990
$ws.Range("J119").Value = 39990
$ws.Range("L119").Value = 39990
$ws.Range("N119").Value = -49666

# Row 132
$ws.Range("H132").Value = 2944.7036
$ws.Range("I132").Value = 1222.8462
$ws.Range("J132").Value = 4543.5713
$ws.Range("K132").Value = 3668.5386
$ws.Range("L132").Value = 13630.7139
$ws.Range("M132").Value = -1138.5386
$ws.Range("N132").Value = -18690.7139

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1954.6897
$ws.Range("I134").Value = 1224.7084
$ws.Range("J134").Value = 5458.6
$ws.Range("K134").Value = 3674.1252
$ws.Range("L134").Value = 16375.8
$ws.Range("M134").Value = -1139.1252
$ws.Range("N134").Value = -21445.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15876264
$ws.Range("I31").Value = 24390864
$ws.Range("J31").Value = 8146.9546
$ws.Range("K31").Value = 24390864
$ws.Range("L31").Value = 8146.9546
$ws.Range("M31").Value = -24390569
$ws.Range("N31").Value = -8736.954600000001

# Row 34
$ws.Range("H34").Value = 15876264
$ws.Range("I34").Value = 24390864
$ws.Range("J34").Value = 8146.9546
$ws.Range("K34").Value = 24390864
$ws.Range("L34").Value = 8146.9546
$ws.Range("M34").Value = -24390662
$ws.Range("N34").Value = -8550.954600000001

# Row 109
$ws.Range("H109").Value = 36866.668
$ws.Range("J109").Value = 36800
$ws.Range("L109").Value = 36800
$ws.Range("N109").Value = -38880

$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 1155.8572
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1155.8572
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3467.5716
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4007.5716

# Row 67
$ws.Range("H67").Value = 1155.8572
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1155.8572
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3467.5716
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5339.571599999999

# Row 68
$ws.Range("H68").Value = 1437244.1
$ws.Range("I68").Value = 2504043.5
$ws.Range("J68").Value = 1168.1538
$ws.Range("K68").Value = 7512130.5
$ws.Range("L68").Value = 3504.4614
$ws.Range("M68").Value = -7511319.5
$ws.Range("N68").Value = -5126.4614

# Row 71
$ws.Range("H71").Value = 1437244.1
$ws.Range("I71").Value = 2504043.5
$ws.Range("J71").Value = 1168.1538
$ws.Range("K71").Value = 22536391.5
$ws.Range("L71").Value = 10513.3842
$ws.Range("M71").Value = -22532335.5
$ws.Range("N71").Value = -18625.3842

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 720.2
$ws.Range("I81").Value = 720.2
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1440.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -379.4000000000001
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 720.2
$ws.Range("I84").Value = 720.2
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7202
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -1898
$ws.Range("N84").ClearContents()

# Row 103
$ws.Range("H103").Value = 22649.6
$ws.Range("J103").Value = 22649.6
$ws.Range("L103").Value = 22649.6
$ws.Range("N103").Value = -24993.6
